# Refresh the cryptos price/volume table with the latest scraped values.
# Note: several "Price" (column D) values look numeric (e.g. "575.32") but
# must stay exact text, matching the original inline-string cells -- a
# leading apostrophe forces Excel to keep them as text instead of coercing
# to a floating-point number (which would lose trailing zeros / introduce
# FP rounding noise, e.g. "2.00" -> 2 or "575.32" -> 575.32000000000005).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.300.89"
$ws.Range("E2").Value = "  +0.70%  "
$ws.Range("D3").Value = "2.770.90"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'575.32"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("D6").Value = "'160.97"
$ws.Range("E6").Value = "  +1.35%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.602"
$ws.Range("E8").Value = "  -1.61%  "
$ws.Range("E9").Value = "  -1.25%  "
$ws.Range("D10").Value = "'5.89"
$ws.Range("E10").Value = "  -1.71%  "
$ws.Range("E11").Value = "  +4.65%  "
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("D13").Value = "3.262.55"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "'27.37"
$ws.Range("E14").Value = "  +1.36%  "
$ws.Range("D15").Value = "63.949.88"
$ws.Range("E15").Value = "  +0.24%  "
$ws.Range("E16").Value = "  -2.14%  "
$ws.Range("D17").Value = "2.781.29"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "'12.24"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("E19").Value = "  -2.05%  "
$ws.Range("D20").Value = "'361.47"
$ws.Range("E20").Value = "  -0.19%  "
$ws.Range("E21").Value = "  -2.63%  "
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("D23").Value = "'0.532"
$ws.Range("E23").Value = "  -6.15%  "
$ws.Range("D24").Value = "'65.20"
$ws.Range("E24").Value = "  -2.11%  "
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("D26").Value = "'8.64"
$ws.Range("E26").Value = "  -0.70%  "
$ws.Range("D27").Value = "'0.998"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Value = "0.0₃0919"
$ws.Range("E28").Value = "  -1.90%  "
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'2.00"
$ws.Range("E29").Value = "  -0.54%  "
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").Value = "'7.39"
$ws.Range("E30").Value = "  +3.98%  "
$ws.Range("E31").Value = "  +10.32%  "
$ws.Range("D32").Value = "'167.46"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("E33").Value = "  +4.38%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "'20.25"
$ws.Range("E35").Value = "  -1.61%  "
$ws.Range("E36").Value = "  -0.06%  "
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("E38").Value = "  -0.05%  "
$ws.Range("D39").Value = "'350.55"
$ws.Range("E39").Value = "  +5.51%  "
$ws.Range("E40").Value = "  +3.90%  "
$ws.Range("E41").Value = "  +0.51%  "
$ws.Range("E42").Value = "  -1.48%  "
$ws.Range("D43").Value = "'22.72"
$ws.Range("E43").Value = "  +3.41%  "
$ws.Range("D44").Value = "'21.71"
$ws.Range("E44").Value = "  -1.54%  "
$ws.Range("D45").Value = "'0.0597"
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("D46").Value = "'137.39"
$ws.Range("E46").Value = "  +0.62%  "
$ws.Range("E47").Value = "  -1.34%  "
$ws.Range("E48").Value = "  -1.20%  "
$ws.Range("E49").Value = "  -0.74%  "
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").Value = "2.139.34"
$ws.Range("E51").Value = "  +0.17%  "
